$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.299.52"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "3.078.61"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("D5").Value = "'522.43"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").Value = "'135.51"
$ws.Range("E6").Value = "  -4.91%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.078.31"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "'0.467"
$ws.Range("E9").Value = "  +4.79%  "
$ws.Range("D10").Value = "'7.28"
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "'0.401"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "3.607.64"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "'25.20"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("D17").Value = "57.371.09"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "3.075.43"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "'5.87"
$ws.Range("E19").Value = "  -4.14%  "
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").Value = "'7.83"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "'349.71"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'68.98"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "0.0₃0864"
$ws.Range("E28").Value = "  -6.92%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").Value = "'5.79"
$ws.Range("E32").Value = "  -9.09%  "
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("D35").Value = "'159.23"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E36").Value = "  -4.46%  "
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("D38").Value = "'25.45"
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("D39").Value = "'1.24"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").Value = "'0.0657"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "'4.05"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("E42").Value = "  -5.68%  "
$ws.Range("D43").Value = "'0.694"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "2.417.47"
$ws.Range("E44").Value = "  +6.24%  "
$ws.Range("D45").Value = "'36.64"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "3.115.72"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'5.97"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'0.943"
$ws.Range("E50").Value = "  -5.79%  "
$ws.Range("D51").Value = "'19.59"
$ws.Range("E51").Value = "  -5.26%  "
